$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.842200040817261
$ws.Range("B1").Value = 3.075647830963135
$ws.Range("C1").Value = 2.800085306167603
$ws.Range("D1").Value = 3.159723997116089
$ws.Range("E1").Value = 2.393796920776367
